# Generate Report for Handoff
# Updates the handoff status for the 6 "Ready for handoff" rows (rows 7,8,9,10,12,14)
# across the Overview, zh-cn and de-de sheets: sets the Priority to "ht" on the
# locale sheets, and refreshes the "Latest Handoff"/"Latest HO Xliff Generate Date"
# timestamps to reflect the newly generated handoff.

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 9, 10, 12, 14)

# Overview sheet: column G = "Latest HO Xliff Generate Date"
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Cells.Item($r, 7).Value = "2016-08-22 11:44:58"
}

# zh-cn sheet: column E = "Priority", column H = "Latest Handoff Datetime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Cells.Item($r, 5).Value = "ht"
    $wsZhCn.Cells.Item($r, 8).Value = "2016-08-22 11:44:53"
}

# de-de sheet: column E = "Priority"
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Cells.Item($r, 5).Value = "ht"
}
